$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.849.62'
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").Value = '2.490.87'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.01%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +2.80%  '
$ws.Range("D9").Value = '2.515.76'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0995'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").Value = '2.937.19'
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").Value = '58.819.59'
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("E17").Value = '  +3.13%  '
$ws.Range("D18").Value = '2.504.68'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.994'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.993'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.84%  '
$ws.Range("D29").Value = '0.0₃0766'
$ws.Range("E29").Value = '  +5.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.68%  '
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("E39").Value = '  +4.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.789'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '280.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("E43").Value = '  +3.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '132.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.594'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0509'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0219'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.18%  '
$ws.Range("D51").Value = '1.759.10'
$ws.Range("E51").Value = '  +2.94%  '
